$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 24
$ws.Cells.Item(24, 1).Value = 24
$ws.Cells.Item(24, 2).Value = "优化一下代码，比如增加注释，对调试信息的管控、改善列表页面判断的顺序，对常量进行全局的管理，去掉一些没有用的代码等等。查找一下隐藏的bug。"
$ws.Cells.Item(24, 4).Value = "完成了调试信息的管控、全局变量进行了管理。其他未完成"
$ws.Cells.Item(24, 5).Value = "ok"
$ws.Cells.Item(24, 6).Value = "2016.5.7"
$ws.Rows.Item(24).RowHeight = 99

# Row 25
$ws.Cells.Item(25, 1).Value = 25
$ws.Cells.Item(25, 2).Value = "第一次进入无法全屏"
$ws.Cells.Item(25, 4).Value = "退出时，保存顶部和进度条隐藏的状态值"
$ws.Cells.Item(25, 5).Value = "ok"
$ws.Cells.Item(25, 6).Value = "2016.5.7"

# Row 26
$ws.Cells.Item(26, 1).Value = 26
$ws.Cells.Item(26, 2).Value = "图片界面一直处于加载界面"
$ws.Cells.Item(26, 4).Value = "判断条件出错，一直无法进入刷新"
$ws.Cells.Item(26, 5).Value = "ok"
$ws.Cells.Item(26, 6).Value = "2016.5.7"

# Row 27
$ws.Cells.Item(27, 1).Value = 27
$ws.Cells.Item(27, 2).Value = "在播放视频的时候去看缩略图片，刷新显示的慢，在音乐界面就很快"
$ws.Cells.Item(27, 3).Value = "估计是内存的占用问题"
$ws.Cells.Item(27, 5).Value = "。。。。"

# Freeze panes at row 2 (split after row 1), keep view scrolled to row 20, select D20
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("A20").Select()
$ws.Range("D20").Select()
